# Penalty Reward System edit: shift the 16-week forecast window by one
# week (each week's date/forecast moves into the previous week's slot,
# and a new W16 week is appended) and refresh the dependent Summary
# sheet figures.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" -----------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$forecastRows = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 232 },
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 149 },
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 153 },
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 142 },
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 147 },
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 142 },
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 140 },
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 140 },
    @{ Row = 10; Date = "2025-03-09"; Forecast = 134 },
    @{ Row = 11; Date = "2025-03-16"; Forecast = 135 },
    @{ Row = 12; Date = "2025-03-23"; Forecast = 134 },
    @{ Row = 13; Date = "2025-03-30"; Forecast = 134 },
    @{ Row = 14; Date = "2025-04-06"; Forecast = 130 },
    @{ Row = 15; Date = "2025-04-13"; Forecast = 126 },
    @{ Row = 16; Date = "2025-04-20"; Forecast = 126 },
    @{ Row = 17; Date = "2025-04-27"; Forecast = 129 }
)

foreach ($item in $forecastRows) {
    $r = $item.Row

    # Column B (Week_Start_Date) is stored as plain text in this sheet,
    # not as a real date serial - force Text formatting before writing so
    # the COM date auto-parser doesn't turn it into a date number.
    $cellB = $wsForecast.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $item.Date

    # Column D (MyForecast) is a genuine numeric column.
    $wsForecast.Cells.Item($r, 4).Value = $item.Forecast
}

# --- Sheet 2: "Summary" ------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$summaryRows = @(
    @{ Row = 2;  Value = "2024-06-23 to 2025-01-05" },
    @{ Row = 6;  Value = "64" },
    @{ Row = 7;  Value = "91" },
    @{ Row = 8;  Value = "2634 units" },
    @{ Row = 9;  Value = "2294" },
    @{ Row = 10; Value = "1246" },
    @{ Row = 11; Value = "676" },
    @{ Row = 12; Value = "232" },
    @{ Row = 13; Value = "2025-01-12" },
    @{ Row = 14; Value = "126" }
)

foreach ($item in $summaryRows) {
    $r = $item.Row
    # Column B on the Summary sheet holds text (dates, counts, "NNNN units"
    # strings, etc.) even when the text looks numeric - force Text so COM
    # doesn't coerce "64"/"2294"/... into real numbers or parse dates.
    $cellB = $wsSummary.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $item.Value
}
